$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores numeric-looking values as plain text (e.g. "582.80").
# Where the refreshed value would otherwise be auto-detected as a number by
# Excel, mark that cell Text first so it keeps the original text representation.

$ws.Range("D2").Value = "62.132.35"
$ws.Range("E2").Value = "  -1.91%  "

$ws.Range("D3").Value = "3.015.83"
$ws.Range("E3").Value = "  -1.55%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.11"
$ws.Range("E5").Value = "  -1.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.94"
$ws.Range("E6").Value = "  -4.35%  "

$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("E8").Value = "  -2.30%  "

$ws.Range("D9").Value = "3.013.87"
$ws.Range("E9").Value = "  -1.42%  "

$ws.Range("E10").Value = "  -4.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.82"
$ws.Range("E11").Value = "  -1.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.443"
$ws.Range("E12").Value = "  -1.61%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000228"
$ws.Range("E13").Value = "  -3.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.52"
$ws.Range("E14").Value = "  -5.40%  "

$ws.Range("E15").Value = "  +2.01%  "

$ws.Range("D16").Value = "3.508.02"
$ws.Range("E16").Value = "  -1.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.03"
$ws.Range("E17").Value = "  -1.86%  "

$ws.Range("D18").Value = "62.060.39"
$ws.Range("E18").Value = "  -1.94%  "

$ws.Range("D19").Value = "3.009.60"
$ws.Range("E19").Value = "  -2.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "464.09"
$ws.Range("E20").Value = "  -3.91%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.92"
$ws.Range("E21").Value = "  -3.67%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.684"
$ws.Range("E22").Value = "  -3.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.46"
$ws.Range("E23").Value = "  -0.74%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.79"
$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.25"
$ws.Range("E25").Value = "  -6.26%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.25"
$ws.Range("E26").Value = "  -4.02%  "

$ws.Range("E27").Value = "  +0.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.92"
$ws.Range("E28").Value = "  -6.33%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.62"
$ws.Range("E30").Value = "  -2.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.16"
$ws.Range("E31").Value = "  -5.26%  "

$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.65"
$ws.Range("E32").Value = "  +5.37%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.08"
$ws.Range("E33").Value = "  -5.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.109"
$ws.Range("E34").Value = "  -2.42%  "

$ws.Range("E35").Value = "  -3.44%  "

$ws.Range("D36").Value = "0.0₃0794"
$ws.Range("E36").Value = "  -3.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.77"
$ws.Range("E37").Value = "  -4.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.11"
$ws.Range("E38").Value = "  -4.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.26"
$ws.Range("E39").Value = "  -0.64%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.04"
$ws.Range("E40").Value = "  -2.46%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.93"
$ws.Range("E41").Value = "  -9.77%  "

$ws.Range("E42").Value = "  -0.29%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.273"
$ws.Range("E43").Value = "  -5.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0353"
$ws.Range("E44").Value = "  -2.45%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "382.55"
$ws.Range("E45").Value = "  -13.04%  "

$ws.Range("D46").Value = "2.755.31"
$ws.Range("E46").Value = "  -2.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "37.64"
$ws.Range("E47").Value = "  -5.87%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.94"
$ws.Range("E48").Value = "  -3.45%  "

$ws.Range("E49").Value = "  +0.08%  "

$ws.Range("E50").Value = "  -1.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.17"
$ws.Range("E51").Value = "  -2.39%  "
